# Update countries & provincias Spain
# Applies the periodic data refresh to the "Pais" sheet:
#  - Updates the "last updated" timestamp in A1
#  - Updates case statistics for several countries (rows 4, 16, 44, 45, 95, 96)
#  - Two pairs of countries swap rank order (Panama/Luxemburgo, Uruguay/Burkina Faso)
#    because their "Casos totales" crossed over; the row whose country moved up
#    receives the freshly-updated figures while the row that moved down keeps the
#    figures the higher-ranked row previously had.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# 1. Refresh timestamp banner
$ws.Range("A1").Value = "Datos actualizados a 13 de Abril de 2020 a las 01:52"

# 2. Estados Unidos (row 4) - figures updated in place
$ws.Range("B4").Value = 560246
$ws.Range("C4").Value = 27367
$ws.Range("D4").Value = 32237
$ws.Range("E4").Value = 505908
$ws.Range("F4").Value = 11766
$ws.Range("G4").Value = 1524
$ws.Range("H4").Value = 22101

# 3. Canada (row 16) - figures updated in place
$ws.Range("B16").Value = 24383
$ws.Range("C16").Value = 1065
$ws.Range("D16").Value = 7172
$ws.Range("E16").Value = 16494
$ws.Range("F16").Value = 557
$ws.Range("G16").Value = 64
$ws.Range("H16").Value = 717

# 4. Panama overtakes Luxemburgo -> swap the two country names and update figures
$ws.Range("A44").Value = "Panama"
$ws.Range("B44").Value = 3400
$ws.Range("C44").Value = 166
$ws.Range("D44").Value = 29
$ws.Range("E44").Value = 3284
$ws.Range("F44").Value = 106
$ws.Range("G44").Value = 8
$ws.Range("H44").Value = 87

$ws.Range("A45").Value = "Luxemburgo"
$ws.Range("B45").Value = 3281
$ws.Range("C45").Value = 11
$ws.Range("D45").Value = 500
$ws.Range("E45").Value = 2715
$ws.Range("F45").Value = 30
$ws.Range("G45").Value = 4
$ws.Range("H45").Value = 66

# 5. Uruguay overtakes Burkina Faso -> swap the two country names and update figures
$ws.Range("A95").Value = "Uruguay"
$ws.Range("B95").Value = 512
$ws.Range("C95").Value = 18
$ws.Range("D95").Value = 231
$ws.Range("E95").Value = 274
$ws.Range("F95").Value = 16
$ws.Range("G95").Value = 0
$ws.Range("H95").Value = 7

$ws.Range("A96").Value = "Burkina Faso"
$ws.Range("B96").Value = 497
$ws.Range("C96").Value = 13
$ws.Range("D96").Value = 161
$ws.Range("E96").Value = 309
$ws.Range("F96").Value = 0
$ws.Range("G96").Value = 0
$ws.Range("H96").Value = 27
